$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 531.3333  # H33: 502.4375 -> 531.3333
$ws.Cells.Item(33, 9).Value = 230.83333  # I33: 218.38461 -> 230.83333
$ws.Cells.Item(33, 11).Value = 230.83333  # K33: 218.38461 -> 230.83333
$ws.Cells.Item(33, 13).Value = -1.833329999999989  # M33: 10.61538999999999 -> -1.833329999999989
$ws.Cells.Item(40, 8).Value = 3325.4211  # H40: 3037.9092 -> 3325.4211
$ws.Cells.Item(40, 9).Value = 0  # I40: 1599 -> 0
$ws.Cells.Item(40, 10).Value = 3325.4211  # J40: 3181.8 -> 3325.4211
$ws.Cells.Item(40, 11).Value = 0  # K40: 1599 -> 0
$ws.Cells.Item(40, 12).Value = 3325.4211  # L40: 3181.8 -> 3325.4211
$ws.Cells.Item(40, 13).ClearContents()  # M40: was -1424
$ws.Cells.Item(40, 14).Value = -3675.4211  # N40: -3531.8 -> -3675.4211
$ws.Cells.Item(104, 8).Value = 196.8  # H104: 221.25 -> 196.8
$ws.Cells.Item(104, 9).Value = 196.8  # I104: 221.25 -> 196.8
$ws.Cells.Item(104, 11).Value = 590.4000000000001  # K104: 663.75 -> 590.4000000000001
$ws.Cells.Item(104, 13).Value = 1156.6  # M104: 1083.25 -> 1156.6
$ws.Cells.Item(107, 8).Value = 634.35486  # H107: 604.6061 -> 634.35486
$ws.Cells.Item(107, 10).Value = 672  # J107: 521 -> 672
$ws.Cells.Item(107, 12).Value = 672  # L107: 521 -> 672
$ws.Cells.Item(107, 14).Value = -4512  # N107: -4361 -> -4512
$ws.Cells.Item(129, 8).Value = 3883.1428  # H129: 3396.8235 -> 3883.1428
$ws.Cells.Item(129, 9).Value = 1198  # I129: 1195.75 -> 1198
$ws.Cells.Item(129, 10).Value = 4330.6665  # J129: 4074.077 -> 4330.6665
$ws.Cells.Item(129, 11).Value = 3594  # K129: 3587.25 -> 3594
$ws.Cells.Item(129, 12).Value = 12991.9995  # L129: 12222.231 -> 12991.9995
$ws.Cells.Item(129, 13).Value = 1406  # M129: 1412.75 -> 1406
$ws.Cells.Item(129, 14).Value = -22991.9995  # N129: -22222.231 -> -22991.9995
$ws.Cells.Item(135, 8).Value = 1839  # H135: 811.75 -> 1839
$ws.Cells.Item(135, 9).Value = 99  # I135: 82.5 -> 99
$ws.Cells.Item(135, 10).Value = 2999  # J135: 2999.5 -> 2999
$ws.Cells.Item(135, 11).Value = 891  # K135: 742.5 -> 891
$ws.Cells.Item(135, 12).Value = 26991  # L135: 26995.5 -> 26991
$ws.Cells.Item(135, 13).Value = 1644  # M135: 1792.5 -> 1644
$ws.Cells.Item(135, 14).Value = -32061  # N135: -32065.5 -> -32061
$ws.Cells.Item(138, 8).Value = 2932.1667  # H138: 2798.923 -> 2932.1667
$ws.Cells.Item(138, 9).Value = 1997.2  # I138: 1864.3334 -> 1997.2
$ws.Cells.Item(138, 11).Value = 5991.6  # K138: 5593.0002 -> 5991.6
$ws.Cells.Item(138, 13).Value = -851.6000000000004  # M138: -453.0002000000004 -> -851.6000000000004
$ws.Cells.Item(141, 8).Value = 5611.25  # H141: 4033.3333 -> 5611.25
$ws.Cells.Item(141, 9).Value = 5611.25  # I141: 4033.3333 -> 5611.25
$ws.Cells.Item(141, 11).Value = 16833.75  # K141: 12099.9999 -> 16833.75
$ws.Cells.Item(141, 13).Value = -11653.75  # M141: -6919.999899999999 -> -11653.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2900  # H61: 2183.1667 -> 2900
$ws.Cells.Item(61, 9).Value = 0  # I61: 1825 -> 0
$ws.Cells.Item(61, 10).Value = 2900  # J61: 2899.5 -> 2900
$ws.Cells.Item(61, 11).Value = 0  # K61: 1825 -> 0
$ws.Cells.Item(61, 12).Value = 2900  # L61: 2899.5 -> 2900
$ws.Cells.Item(61, 13).ClearContents()  # M61: was -1613
$ws.Cells.Item(61, 14).Value = -3324  # N61: -3323.5 -> -3324
$ws.Cells.Item(74, 8).Value = 1966  # H74: 1822.1111 -> 1966
$ws.Cells.Item(74, 9).Value = 1699.8572  # I74: 1687.375 -> 1699.8572
$ws.Cells.Item(74, 10).Value = 2897.5  # J74: 2900 -> 2897.5
$ws.Cells.Item(74, 11).Value = 1699.8572  # K74: 1687.375 -> 1699.8572
$ws.Cells.Item(74, 12).Value = 2897.5  # L74: 2900 -> 2897.5
$ws.Cells.Item(74, 13).Value = -825.8571999999999  # M74: -813.375 -> -825.8571999999999
$ws.Cells.Item(74, 14).Value = -4645.5  # N74: -4648 -> -4645.5
$ws.Cells.Item(77, 8).Value = 1966  # H77: 1822.1111 -> 1966
$ws.Cells.Item(77, 9).Value = 1699.8572  # I77: 1687.375 -> 1699.8572
$ws.Cells.Item(77, 10).Value = 2897.5  # J77: 2900 -> 2897.5
$ws.Cells.Item(77, 11).Value = 8499.286  # K77: 8436.875 -> 8499.286
$ws.Cells.Item(77, 12).Value = 14487.5  # L77: 14500 -> 14487.5
$ws.Cells.Item(77, 13).Value = -4131.286  # M77: -4068.875 -> -4131.286
$ws.Cells.Item(77, 14).Value = -23223.5  # N77: -23236 -> -23223.5
$ws.Cells.Item(136, 8).Value = 2900  # H136: 2183.1667 -> 2900
$ws.Cells.Item(136, 9).Value = 0  # I136: 1825 -> 0
$ws.Cells.Item(136, 10).Value = 2900  # J136: 2899.5 -> 2900
$ws.Cells.Item(136, 11).Value = 0  # K136: 5475 -> 0
$ws.Cells.Item(136, 12).Value = 8700  # L136: 8698.5 -> 8700
$ws.Cells.Item(136, 13).ClearContents()  # M136: was -2925
$ws.Cells.Item(136, 14).Value = -13800  # N136: -13798.5 -> -13800

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2173.4285  # H86: 2331.4285 -> 2173.4285
$ws.Cells.Item(86, 9).Value = 2293.2  # I86: 2504.4 -> 2293.2
$ws.Cells.Item(86, 10).Value = 1874  # J86: 1899 -> 1874
$ws.Cells.Item(86, 11).Value = 2293.2  # K86: 2504.4 -> 2293.2
$ws.Cells.Item(86, 12).Value = 1874  # L86: 1899 -> 1874
$ws.Cells.Item(86, 13).Value = -1170.2  # M86: -1381.4 -> -1170.2
$ws.Cells.Item(86, 14).Value = -4120  # N86: -4145 -> -4120
$ws.Cells.Item(89, 8).Value = 2173.4285  # H89: 2331.4285 -> 2173.4285
$ws.Cells.Item(89, 9).Value = 2293.2  # I89: 2504.4 -> 2293.2
$ws.Cells.Item(89, 10).Value = 1874  # J89: 1899 -> 1874
$ws.Cells.Item(89, 11).Value = 11466  # K89: 12522 -> 11466
$ws.Cells.Item(89, 12).Value = 9370  # L89: 9495 -> 9370
$ws.Cells.Item(89, 13).Value = -5850  # M89: -6906 -> -5850
$ws.Cells.Item(89, 14).Value = -20602  # N89: -20727 -> -20602
$ws.Cells.Item(107, 8).Value = 753.6  # H107: 770.5333000000001 -> 753.6
$ws.Cells.Item(107, 9).Value = 354.58334  # I107: 382.63635 -> 354.58334
$ws.Cells.Item(107, 10).Value = 2349.6667  # J107: 1837.25 -> 2349.6667
$ws.Cells.Item(107, 11).Value = 354.58334  # K107: 382.63635 -> 354.58334
$ws.Cells.Item(107, 12).Value = 2349.6667  # L107: 1837.25 -> 2349.6667
$ws.Cells.Item(107, 13).Value = 1565.41666  # M107: 1537.36365 -> 1565.41666
$ws.Cells.Item(107, 14).Value = -6189.6667  # N107: -5677.25 -> -6189.6667
$ws.Cells.Item(134, 8).Value = 4476.222  # H134: 4330.0527 -> 4476.222
$ws.Cells.Item(134, 9).Value = 4363.0586  # I134: 4215.0557 -> 4363.0586
$ws.Cells.Item(134, 11).Value = 13089.1758  # K134: 12645.1671 -> 13089.1758
$ws.Cells.Item(134, 13).Value = -10554.1758  # M134: -10110.1671 -> -10554.1758

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 537.5  # H22: 616.6667 -> 537.5
$ws.Cells.Item(22, 9).Value = 537.5  # I22: 650 -> 537.5
$ws.Cells.Item(22, 10).Value = 0  # J22: 600 -> 0
$ws.Cells.Item(22, 11).Value = 537.5  # K22: 650 -> 537.5
$ws.Cells.Item(22, 12).Value = 0  # L22: 600 -> 0
$ws.Cells.Item(22, 13).Value = -187.5  # M22: -300 -> -187.5
$ws.Cells.Item(22, 14).ClearContents()  # N22: was -1300
$ws.Cells.Item(31, 8).Value = 3586  # H31: 4117.467 -> 3586
$ws.Cells.Item(31, 10).Value = 3304.923  # J31: 3978.5454 -> 3304.923
$ws.Cells.Item(31, 12).Value = 3304.923  # L31: 3978.5454 -> 3304.923
$ws.Cells.Item(31, 14).Value = -3894.923  # N31: -4568.5454 -> -3894.923
$ws.Cells.Item(34, 8).Value = 3586  # H34: 4117.467 -> 3586
$ws.Cells.Item(34, 10).Value = 3304.923  # J34: 3978.5454 -> 3304.923
$ws.Cells.Item(34, 12).Value = 3304.923  # L34: 3978.5454 -> 3304.923
$ws.Cells.Item(34, 14).Value = -3708.923  # N34: -4382.5454 -> -3708.923
$ws.Cells.Item(58, 8).Value = 4025  # H58: 3639 -> 4025
$ws.Cells.Item(58, 9).Value = 2550  # I58: 2398.3333 -> 2550
$ws.Cells.Item(58, 11).Value = 2550  # K58: 2398.3333 -> 2550
$ws.Cells.Item(58, 13).Value = -2347  # M58: -2195.3333 -> -2347
$ws.Cells.Item(99, 8).Value = 5500  # H99: 8000 -> 5500
$ws.Cells.Item(99, 9).Value = 3000  # I99: 0 -> 3000
$ws.Cells.Item(99, 11).Value = 3000  # K99: 0 -> 3000
$ws.Cells.Item(99, 13).Value = -1502  # M99: None -> -1502
$ws.Cells.Item(107, 8).Value = 549.7778  # H107: 550 -> 549.7778
$ws.Cells.Item(107, 9).Value = 395.8  # I107: 375.25 -> 395.8
$ws.Cells.Item(107, 10).Value = 742.25  # J107: 689.8 -> 742.25
$ws.Cells.Item(107, 11).Value = 395.8  # K107: 375.25 -> 395.8
$ws.Cells.Item(107, 12).Value = 742.25  # L107: 689.8 -> 742.25
$ws.Cells.Item(107, 13).Value = 1524.2  # M107: 1544.75 -> 1524.2
$ws.Cells.Item(107, 14).Value = -4582.25  # N107: -4529.8 -> -4582.25
$ws.Cells.Item(126, 8).Value = 5500  # H126: 8000 -> 5500
$ws.Cells.Item(126, 9).Value = 3000  # I126: 0 -> 3000
$ws.Cells.Item(126, 11).Value = 9000  # K126: 0 -> 9000
$ws.Cells.Item(126, 13).Value = -6530  # M126: None -> -6530
$ws.Cells.Item(132, 8).Value = 3261.4  # H132: 4866.6665 -> 3261.4
$ws.Cells.Item(132, 9).Value = 3261.4  # I132: 4866.6665 -> 3261.4
$ws.Cells.Item(132, 11).Value = 9784.200000000001  # K132: 14599.9995 -> 9784.200000000001
$ws.Cells.Item(132, 13).Value = -7254.200000000001  # M132: -12069.9995 -> -7254.200000000001
$ws.Cells.Item(134, 8).Value = 3113.1  # H134: 3932.1538 -> 3113.1
$ws.Cells.Item(134, 9).Value = 3015  # I134: 2813.6 -> 3015
$ws.Cells.Item(134, 10).Value = 3996  # J134: 7660.6665 -> 3996
$ws.Cells.Item(134, 11).Value = 9045  # K134: 8440.799999999999 -> 9045
$ws.Cells.Item(134, 12).Value = 11988  # L134: 22981.9995 -> 11988
$ws.Cells.Item(134, 13).Value = -6510  # M134: -5905.799999999999 -> -6510
$ws.Cells.Item(134, 14).Value = -17058  # N134: -28051.9995 -> -17058
$ws.Cells.Item(136, 8).Value = 4025  # H136: 3639 -> 4025
$ws.Cells.Item(136, 9).Value = 2550  # I136: 2398.3333 -> 2550
$ws.Cells.Item(136, 11).Value = 7650  # K136: 7194.999899999999 -> 7650
$ws.Cells.Item(136, 13).Value = -5100  # M136: -4644.999899999999 -> -5100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 147  # H33: 133.75 -> 147
$ws.Cells.Item(33, 10).Value = 200  # J33: 0 -> 200
$ws.Cells.Item(33, 12).Value = 1200  # L33: 0 -> 1200
$ws.Cells.Item(33, 14).Value = -1766  # N33: None -> -1766
$ws.Cells.Item(97, 8).Value = 1209.8572  # H97: 833.5454999999999 -> 1209.8572
$ws.Cells.Item(97, 10).Value = 479.66666  # J97: 305.57144 -> 479.66666
$ws.Cells.Item(97, 12).Value = 1438.99998  # L97: 916.71432 -> 1438.99998
$ws.Cells.Item(97, 14).Value = -2430.99998  # N97: -1908.71432 -> -2430.99998
$ws.Cells.Item(117, 8).Value = 730.7143  # H117: 727.5 -> 730.7143
$ws.Cells.Item(117, 10).Value = 763.5  # J117: 777 -> 763.5
$ws.Cells.Item(117, 12).Value = 2290.5  # L117: 2331 -> 2290.5
$ws.Cells.Item(117, 14).Value = -9174.5  # N117: -9215 -> -9174.5
$ws.Cells.Item(120, 8).Value = 0  # H120: 3000 -> 0
$ws.Cells.Item(120, 9).Value = 0  # I120: 3000 -> 0
$ws.Cells.Item(120, 11).Value = 0  # K120: 9000 -> 0
$ws.Cells.Item(120, 13).ClearContents()  # M120: was -4162
$ws.Cells.Item(140, 8).Value = 843.75  # H140: 2000 -> 843.75
$ws.Cells.Item(140, 9).Value = 843.75  # I140: 2000 -> 843.75
$ws.Cells.Item(140, 11).Value = 2531.25  # K140: 6000 -> 2531.25
$ws.Cells.Item(140, 13).Value = 2648.75  # M140: -820 -> 2648.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7801.75  # H70: 7869 -> 7801.75
$ws.Cells.Item(70, 9).Value = 7801.75  # I70: 7869 -> 7801.75
$ws.Cells.Item(70, 11).Value = 7801.75  # K70: 7869 -> 7801.75
$ws.Cells.Item(70, 13).Value = -7531.75  # M70: -7599 -> -7531.75
$ws.Cells.Item(73, 8).Value = 7801.75  # H73: 7869 -> 7801.75
$ws.Cells.Item(73, 9).Value = 7801.75  # I73: 7869 -> 7801.75
$ws.Cells.Item(73, 11).Value = 7801.75  # K73: 7869 -> 7801.75
$ws.Cells.Item(73, 13).Value = -6865.75  # M73: -6933 -> -6865.75
$ws.Cells.Item(80, 8).Value = 7299.8  # H80: 6916.3335 -> 7299.8
$ws.Cells.Item(80, 10).Value = 8000  # J80: 7249.75 -> 8000
$ws.Cells.Item(80, 12).Value = 8000  # L80: 7249.75 -> 8000
$ws.Cells.Item(80, 14).Value = -9996  # N80: -9245.75 -> -9996
$ws.Cells.Item(83, 8).Value = 7299.8  # H83: 6916.3335 -> 7299.8
$ws.Cells.Item(83, 10).Value = 8000  # J83: 7249.75 -> 8000
$ws.Cells.Item(83, 12).Value = 40000  # L83: 36248.75 -> 40000
$ws.Cells.Item(83, 14).Value = -49984  # N83: -46232.75 -> -49984
$ws.Cells.Item(132, 8).Value = 1012  # H132: 0 -> 1012
$ws.Cells.Item(132, 9).Value = 1012  # I132: 0 -> 1012
$ws.Cells.Item(132, 11).Value = 3036  # K132: 0 -> 3036
$ws.Cells.Item(132, 13).Value = -506  # M132: None -> -506

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2904.125  # H16: 3841.625 -> 2904.125
$ws.Cells.Item(16, 9).Value = 3455.5  # I16: 4176.143 -> 3455.5
$ws.Cells.Item(16, 10).Value = 1250  # J16: 1500 -> 1250
$ws.Cells.Item(16, 11).Value = 3455.5  # K16: 4176.143 -> 3455.5
$ws.Cells.Item(16, 12).Value = 1250  # L16: 1500 -> 1250
$ws.Cells.Item(16, 13).Value = -3285.5  # M16: -4006.143 -> -3285.5
$ws.Cells.Item(16, 14).Value = -1590  # N16: -1840 -> -1590
$ws.Cells.Item(40, 8).Value = 6032.8335  # H40: 6326.909 -> 6032.8335
$ws.Cells.Item(40, 9).Value = 4710.5557  # I40: 4949.625 -> 4710.5557
$ws.Cells.Item(40, 11).Value = 4710.5557  # K40: 4949.625 -> 4710.5557
$ws.Cells.Item(40, 13).Value = -4574.5557  # M40: -4813.625 -> -4574.5557
$ws.Cells.Item(136, 8).Value = 3602.6  # H136: 3669.3333 -> 3602.6
$ws.Cells.Item(136, 9).Value = 3602.6  # I136: 3669.3333 -> 3602.6
$ws.Cells.Item(136, 11).Value = 10807.8  # K136: 11007.9999 -> 10807.8
$ws.Cells.Item(136, 13).Value = -8257.799999999999  # M136: -8457.999899999999 -> -8257.799999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 0  # H81: 1997.5 -> 0
$ws.Cells.Item(81, 9).Value = 0  # I81: 1997.5 -> 0
$ws.Cells.Item(81, 11).Value = 0  # K81: 3995 -> 0
$ws.Cells.Item(81, 13).ClearContents()  # M81: was -2934
$ws.Cells.Item(84, 8).Value = 0  # H84: 1997.5 -> 0
$ws.Cells.Item(84, 9).Value = 0  # I84: 1997.5 -> 0
$ws.Cells.Item(84, 11).Value = 0  # K84: 19975 -> 0
$ws.Cells.Item(84, 13).ClearContents()  # M84: was -14671
$ws.Cells.Item(136, 8).Value = 2998.5  # H136: 2500.25 -> 2998.5
$ws.Cells.Item(136, 9).Value = 2998.5  # I136: 2000.3334 -> 2998.5
$ws.Cells.Item(136, 10).Value = 0  # J136: 4000 -> 0
$ws.Cells.Item(136, 11).Value = 8995.5  # K136: 6001.0002 -> 8995.5
$ws.Cells.Item(136, 12).Value = 0  # L136: 12000 -> 0
$ws.Cells.Item(136, 13).Value = -6445.5  # M136: -3451.0002 -> -6445.5
$ws.Cells.Item(136, 14).ClearContents()  # N136: was -17100
